$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. "总计" (summary) sheet: insert a new row for 2022-Q3 right after the
#    header row, shifting all the existing quarterly rows down by one.
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)
$ws1.Rows.Item(2).Insert()

# The freshly inserted row inherited formatting from the row above (the
# header); strip that so the new data row matches the plain data rows,
# then re-apply just the index-column style copied from an existing data
# row (column A carries a distinct style in this sheet).
$ws1.Range("A2:D2").ClearFormats()
$ws1.Range("A3").Copy()
$ws1.Range("A2").PasteSpecial(-4122)

$ws1.Range("A2").Value = 0
$ws1.Range("B2").Value = "2022-Q3"
$ws1.Range("C2").Value = 4
$ws1.Range("D2").Value = 0.28

# Re-number the index column (A) for the rows that shifted down.
for ($r = 3; $r -le 9; $r++) {
    $ws1.Cells.Item($r, 1).Value = $r - 2
}

# ---------------------------------------------------------------------
# 2. Insert a brand-new "2022-Q3" worksheet right before "2022-Q2" and
#    populate it with the quarter's fund holdings.
# ---------------------------------------------------------------------
$sheetBefore = $wb.Worksheets.Item("2022-Q2")
$newSheet = $wb.Worksheets.Add($sheetBefore)
$newSheet.Name = "2022-Q3"
$newSheet.Outline.SummaryRow = 1
$newSheet.Outline.SummaryColumn = 1

# Match the look of the other quarterly sheets: copy the header-row style
# and the index-column style from the neighbouring "2022-Q2" sheet.
$srcSheet = $wb.Worksheets.Item("2022-Q2")
$srcSheet.Range("A1:H1").Copy()
$newSheet.Range("A1:H1").PasteSpecial(-4122)
$srcSheet.Range("A2").Copy()
$newSheet.Range("A2:A5").PasteSpecial(-4122)

$newSheet.Range("B1").Value = "基金代码"
$newSheet.Range("C1").Value = "基金名称"
$newSheet.Range("D1").Value = "基金规模"
$newSheet.Range("E1").Value = "股票总仓位"
$newSheet.Range("F1").Value = "仓位占比"
$newSheet.Range("G1").Value = "持有市值(亿元)"
$newSheet.Range("H1").Value = "仓位排名"

$newSheet.Range("A2").Value = 0
$newSheet.Range("B2").Value = "'501029"
$newSheet.Range("C2").Value = "华宝标普中国A股红利机会指数（LOF）A"
$newSheet.Range("D2").Value = "'8.11"
$newSheet.Range("E2").Value = "'94.26"
$newSheet.Range("F2").Value = "'2.14"
$newSheet.Range("G2").Value = "'0.1736"
$newSheet.Range("H2").Value = 1

$newSheet.Range("A3").Value = 1
$newSheet.Range("B3").Value = "'005125"
$newSheet.Range("C3").Value = "华宝标普中国A股红利机会指数C"
$newSheet.Range("D3").Value = "'3.38"
$newSheet.Range("E3").Value = "'94.26"
$newSheet.Range("F3").Value = "'2.14"
$newSheet.Range("G3").Value = "'0.0723"
$newSheet.Range("H3").Value = 1

$newSheet.Range("A4").Value = 2
$newSheet.Range("B4").Value = "'003845"
$newSheet.Range("C4").Value = "汇安丰恒灵活配置混合A"
$newSheet.Range("D4").Value = "'1.27"
$newSheet.Range("E4").Value = "'61.83"
$newSheet.Range("F4").Value = "'2.72"
$newSheet.Range("G4").Value = "'0.0345"
$newSheet.Range("H4").Value = 9

$newSheet.Range("A5").Value = 3
$newSheet.Range("B5").Value = "'003846"
$newSheet.Range("C5").Value = "汇安丰恒灵活配置混合C"
$newSheet.Range("D5").Value = "'0.00"
$newSheet.Range("E5").Value = "'61.83"
$newSheet.Range("F5").Value = "'2.72"
$newSheet.Range("G5").Value = 0
$newSheet.Range("H5").Value = 9

# The leading apostrophes above (needed so numeric-looking strings like
# "501029"/"8.11" are stored as text, matching the source data) leave a
# "quote prefix" flag on the cell style. Reset those cells back to the
# plain "Normal" style so they look like the rest of the (unstyled) data
# cells instead of picking up a stray custom style.
$newSheet.Range("B2:B5").Style = "Normal"
$newSheet.Range("D2:G5").Style = "Normal"
